$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the password value "T0SMKTH*" -> "T0SMKTH#" across the whole
# B2:B247 password column (all rows share this single value).
$ws.Range("B2:B247").Value = "T0SMKTH#"

# Restore the scrolled/selected view: no frozen topLeftCell scroll and the
# active selection is now the full password column B2:B247.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B2:B247").Select()
